$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4037090.8
$ws.Range("I62").Value = 5685355.5
$ws.Range("J62").Value = 7999.8887
$ws.Range("K62").Value = 5685355.5
$ws.Range("L62").Value = 7999.8887
$ws.Range("M62").Value = -5684731.5
$ws.Range("N62").Value = -9247.8887
$ws.Range("H64").Value = 6167.5586
$ws.Range("J64").Value = 6700
$ws.Range("L64").Value = 6700
$ws.Range("N64").Value = -7196
$ws.Range("H65").Value = 4037090.8
$ws.Range("I65").Value = 5685355.5
$ws.Range("J65").Value = 7999.8887
$ws.Range("K65").Value = 28426777.5
$ws.Range("L65").Value = 39999.4435
$ws.Range("M65").Value = -28423657.5
$ws.Range("N65").Value = -46239.4435
$ws.Range("H67").Value = 6167.5586
$ws.Range("J67").Value = 6700
$ws.Range("L67").Value = 6700
$ws.Range("N67").Value = -8416
$ws.Range("H80").Value = 819.36
$ws.Range("I80").Value = 1036.1818
$ws.Range("J80").Value = 649
$ws.Range("K80").Value = 3108.5454
$ws.Range("L80").Value = 1947
$ws.Range("M80").Value = -2110.5454
$ws.Range("N80").Value = -3943
$ws.Range("H83").Value = 819.36
$ws.Range("I83").Value = 1036.1818
$ws.Range("J83").Value = 649
$ws.Range("K83").Value = 9325.636200000001
$ws.Range("L83").Value = 5841
$ws.Range("M83").Value = -4333.636200000001
$ws.Range("N83").Value = -15825
$ws.Range("H112").Value = 1710.5555
$ws.Range("J112").Value = 1710.5555
$ws.Range("L112").Value = 5131.666499999999
$ws.Range("N112").Value = -7347.666499999999
$ws.Range("H113").Value = 2860.5557
$ws.Range("I113").Value = 2860.5557
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2860.5557
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 393.4443000000001
$ws.Range("N113").ClearContents()  # was -10174.5
$ws.Range("H127").Value = 2449.0908
$ws.Range("I127").Value = 2095.4285
$ws.Range("J127").Value = 3068
$ws.Range("K127").Value = 6286.2855
$ws.Range("L127").Value = 9204
$ws.Range("M127").Value = -1326.2855
$ws.Range("N127").Value = -19124
$ws.Range("H132").Value = 2211
$ws.Range("I132").Value = 2240.3333
$ws.Range("K132").Value = 6720.999899999999
$ws.Range("M132").Value = -4190.999899999999
$ws.Range("H136").Value = 49874.6
$ws.Range("J136").Value = 49874.6
$ws.Range("L136").Value = 49874.6
$ws.Range("N136").Value = -60074.6
$ws.Range("H137").Value = 1548.9131
$ws.Range("I137").Value = 1247.6123
$ws.Range("K137").Value = 3742.8369
$ws.Range("M137").Value = -1192.8369
$ws.Range("H138").Value = 5947.3735
$ws.Range("J138").Value = 7341.46
$ws.Range("L138").Value = 22024.38
$ws.Range("N138").Value = -32304.38

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2525.724
$ws.Range("I122").Value = 1038.9166
$ws.Range("J122").Value = 4958.6816
$ws.Range("K122").Value = 3116.7498
$ws.Range("L122").Value = 14876.0448
$ws.Range("M122").Value = -666.7498000000001
$ws.Range("N122").Value = -19776.0448

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 111000
$ws.Range("J59").Value = 111000
$ws.Range("L59").Value = 111000
$ws.Range("N59").Value = -112694
$ws.Range("H105").Value = 3787.375
$ws.Range("I105").Value = 4070.5715
$ws.Range("K105").Value = 4070.5715
$ws.Range("M105").Value = -2323.5715
$ws.Range("H107").Value = 911789
$ws.Range("I107").Value = 2185.111
$ws.Range("K107").Value = 2185.111
$ws.Range("M107").Value = -265.1109999999999
$ws.Range("H123").Value = 49997.75
$ws.Range("J123").Value = 49997.75
$ws.Range("L123").Value = 49997.75
$ws.Range("N123").Value = -59797.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 254699.75
$ws.Range("I31").Value = 1900
$ws.Range("J31").Value = 507499.5
$ws.Range("K31").Value = 1900
$ws.Range("L31").Value = 507499.5
$ws.Range("M31").Value = -1605
$ws.Range("N31").Value = -508089.5
$ws.Range("H34").Value = 254699.75
$ws.Range("I34").Value = 1900
$ws.Range("J34").Value = 507499.5
$ws.Range("K34").Value = 1900
$ws.Range("L34").Value = 507499.5
$ws.Range("M34").Value = -1698
$ws.Range("N34").Value = -507903.5
$ws.Range("H41").Value = 26223.4
$ws.Range("I41").Value = 3705.6667
$ws.Range("K41").Value = 3705.6667
$ws.Range("M41").Value = -3277.6667
$ws.Range("H99").Value = 7135.727
$ws.Range("I99").Value = 6415.5
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 6415.5
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = -4917.5
$ws.Range("N99").Value = -10996
$ws.Range("H126").Value = 7135.727
$ws.Range("I126").Value = 6415.5
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 19246.5
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -16776.5
$ws.Range("N126").Value = -28940
$ws.Range("H139").Value = 98705.8
$ws.Range("J139").Value = 98705.8
$ws.Range("L139").Value = 98705.8
$ws.Range("N139").Value = -108985.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 56
$ws.Range("I10").Value = 56
$ws.Range("K10").Value = 168
$ws.Range("M10").Value = -29
$ws.Range("H13").Value = 4379.8
$ws.Range("I13").Value = 2000
$ws.Range("K13").Value = 6000
$ws.Range("M13").Value = -5832
$ws.Range("H113").Value = 1454268.9
$ws.Range("I113").Value = 3528443
$ws.Range("J113").Value = 2347.0667
$ws.Range("K113").Value = 10585329
$ws.Range("L113").Value = 7041.2001
$ws.Range("M113").Value = -10583159
$ws.Range("N113").Value = -11381.2001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 27512542
$ws.Range("I62").Value = 55000000
$ws.Range("J62").Value = 25084
$ws.Range("K62").Value = 55000000
$ws.Range("L62").Value = 25084
$ws.Range("M62").Value = -54999314
$ws.Range("N62").Value = -26456
$ws.Range("H65").Value = 27512542
$ws.Range("I65").Value = 55000000
$ws.Range("J65").Value = 25084
$ws.Range("K65").Value = 165000000
$ws.Range("L65").Value = 75252
$ws.Range("M65").Value = -164996568
$ws.Range("N65").Value = -82116
$ws.Range("H80").Value = 628795.0600000001
$ws.Range("I80").Value = 718177.6
$ws.Range("K80").Value = 718177.6
$ws.Range("M80").Value = -717179.6
$ws.Range("H83").Value = 628795.0600000001
$ws.Range("I83").Value = 718177.6
$ws.Range("K83").Value = 3590888
$ws.Range("M83").Value = -3585896
$ws.Range("H97").Value = 681.58826
$ws.Range("I97").Value = 732.8
$ws.Range("J97").Value = 297.5
$ws.Range("K97").Value = 732.8
$ws.Range("L97").Value = 297.5
$ws.Range("M97").Value = -236.8
$ws.Range("N97").Value = -1289.5
$ws.Range("H102").Value = 2914.0715
$ws.Range("I102").Value = 1568.5883
$ws.Range("J102").Value = 4993.4546
$ws.Range("K102").Value = 1568.5883
$ws.Range("L102").Value = 4993.4546
$ws.Range("M102").Value = 53.41170000000011
$ws.Range("N102").Value = -8237.454600000001
$ws.Range("H126").Value = 3793.7058
$ws.Range("I126").Value = 3332.8333
$ws.Range("K126").Value = 9998.499899999999
$ws.Range("M126").Value = -7528.499899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5074.6904
$ws.Range("I122").Value = 4494.3076
$ws.Range("J122").Value = 6017.8125
$ws.Range("K122").Value = 13482.9228
$ws.Range("L122").Value = 18053.4375
$ws.Range("M122").Value = -11032.9228
$ws.Range("N122").Value = -22953.4375
$ws.Range("H137").Value = 54996.668
$ws.Range("I137").Value = 57495
$ws.Range("K137").Value = 57495
$ws.Range("M137").Value = -52395
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()  # was -110270

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 128277
$ws.Range("J133").Value = 128277
$ws.Range("L133").Value = 128277
$ws.Range("N133").Value = -138397
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = 0
$ws.Range("H139").Value = 54994.668
$ws.Range("J139").Value = 54994.668
$ws.Range("L139").Value = 54994.668
$ws.Range("N139").Value = -65274.668
